$d = $word.ActiveDocument

$replacements = @(
    @("PRINCIPAL CONSULTANT - Siege Analytics, Washington, DC | January 2014 – Present",
      "PARTNER - Siege Analytics, Washington, DC | January 2014 – Present"),
    @("PRINCIPAL CONSULTANT - Clarity and Rigour, Washington, DC | 2012 – 2014",
      "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014"),
    @("SENIOR CONSULTANT - Helm, Washington, DC | 2010 – 2012",
      "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012"),
    @("CONSULTANT - GSD&M, Austin, TX | 2008 – 2010",
      "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010"),
    @("SENIOR CONSULTANT - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008",
      "RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008"),
    @("CONSULTANT - Salsa Labs, Inc., Washington, DC | 2004 – 2006",
      "SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006"),
    @("CONSULTANT - The Praxis Project, Oakland, CA | 2002 – 2004",
      "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004"),
    @("ANALYST - Lake Research Partners, Washington, DC | 2001 – 2002",
      "PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002"),
    @("FIELD COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001",
      "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
